$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update column F values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 1711
$ws1.Range("F11").Value = 1768
$ws1.Range("F13").Value = 108
$ws1.Range("F14").Value = 421
$ws1.Range("F21").Value = 735
$ws1.Range("F22").Value = 306

# Sheet "全部类型" (fourth sheet) - update column F values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 1711
$ws4.Range("F12").Value = 1768
$ws4.Range("F14").Value = 108
$ws4.Range("F15").Value = 421
$ws4.Range("F22").Value = 735
$ws4.Range("F23").Value = 306
